# NumerosSelecionados.xlsx - apply commit:
# "Alteracoes e Correcoes feitas, faltando estetica de apagar Escolha 10 numeros"
#
# - Fill in previously-empty idPagamento (column D) values for rows 29-33.
# - Complete row 34 (add idPagamento + the "Escolha 10 numeros" / Pagamento columns).
# - Append new rows 35-46 with the same shape of data.
# - Dimension / used range grows from A1:O34 to A1:O46 automatically once O46 is set.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns E..N hold plain numbers ("Escolha 10 numeros"); A/C/D/O hold text.
$numCols = @("E","F","G","H","I","J","K","L","M","N")

function Set-TextCell($sheet, $col, $row, $val) {
    $addr = "$col$row"
    $sheet.Range($addr).NumberFormat = "@"
    if ($val -eq "") {
        # Force an explicit empty *text* cell (t="str") instead of an empty/blank cell.
        $sheet.Range($addr).Formula = '=""'
    } else {
        $sheet.Range($addr).Value2 = $val
    }
}

function Set-NumberCell($sheet, $col, $row, $val) {
    $sheet.Range("$col$row").Value2 = $val
}

function Set-DataRow($sheet, $row, $data) {
    Set-TextCell   $sheet "A" $row $data["A"]
    Set-NumberCell $sheet "B" $row $data["B"]
    Set-TextCell   $sheet "C" $row $data["C"]
    Set-TextCell   $sheet "D" $row $data["D"]
    foreach ($col in $numCols) {
        Set-NumberCell $sheet $col $row $data[$col]
    }
    Set-TextCell $sheet "O" $row $data["O"]
}

# --- Rows 29-33: only idPagamento (column D) was previously blank/empty and now gets a value ---
$dUpdates = @{
    29 = "76999025539"
    30 = "77243067266"
    31 = "77243137454"
    32 = "77000631317"
    33 = "77001090761"
}
foreach ($r in $dUpdates.Keys) {
    Set-TextCell $ws "D" $r $dUpdates[$r]
}

# --- Row 34 (finish it) and brand-new rows 35-46 ---
$newRows = @(
    @{ Row=34; A="Vitor Ito"; B=1578424633; C="11966554411"; D="77001235025"; E=1;  F=2;  G=3;  H=4;  I=5;  J=6;  K=7;  L=8; M=9; N=10; O="Não" },
    @{ Row=35; A="Vitor Ito"; B=1578424633; C="11988776655"; D="77001321047"; E=1;  F=2;  G=3;  H=4;  I=5;  J=6;  K=7;  L=8; M=9; N=10; O="Não" },
    @{ Row=36; A="Vitor Ito"; B=1578424633; C="11966554411"; D="77243935384"; E=18; F=29; G=31; H=32; I=34; J=44; K=45; L=47; M=48; N=50; O="Não" },
    @{ Row=37; A="Vitor Ito"; B=1578424633; C="12965487563"; D="77001321967"; E=1;  F=2;  G=3;  H=4;  I=5;  J=7;  K=8;  L=9; M=10; N=12; O="Não" },
    @{ Row=38; A="Vitor Ito"; B=1578424633; C="11955448877"; D="77244329042"; E=1;  F=2;  G=3;  H=4;  I=5;  J=6;  K=7;  L=8; M=9; N=10; O="Não" },
    @{ Row=39; A="Vitor Ito"; B=1578424633; C="11966554422"; D="77001940671"; E=1;  F=2;  G=3;  H=4;  I=5;  J=6;  K=7;  L=8; M=9; N=10; O="Não" },
    @{ Row=40; A="Vitor Ito"; B=1578424633; C="11988775566"; D="77244499150"; E=1;  F=2;  G=3;  H=4;  I=5;  J=6;  K=7;  L=8; M=9; N=10; O="Não" },
    @{ Row=41; A="Vitor Ito"; B=1578424633; C="11966554400"; D="";            E=1;  F=2;  G=3;  H=4;  I=5;  J=6;  K=7;  L=8; M=9; N=10; O="Não" },
    @{ Row=42; A="Vitor Ito"; B=1578424633; C="";            D="";            E=1;  F=2;  G=3;  H=4;  I=5;  J=6;  K=7;  L=8; M=9; N=10; O="Não" },
    @{ Row=43; A="Vitor Ito"; B=1578424633; C="11966554422"; D="";            E=1;  F=2;  G=3;  H=4;  I=5;  J=6;  K=7;  L=8; M=9; N=10; O="Não" },
    @{ Row=44; A="Vitor Ito"; B=1578424633; C="";            D="";            E=1;  F=2;  G=3;  H=4;  I=5;  J=6;  K=7;  L=8; M=9; N=10; O="Não" },
    @{ Row=45; A="Vitor Ito"; B=1578424633; C="11988553322"; D="";            E=1;  F=2;  G=3;  H=4;  I=5;  J=6;  K=7;  L=8; M=9; N=10; O="Não" },
    @{ Row=46; A="Vitor Ito"; B=1578424633; C="11966548087"; D="";            E=1;  F=2;  G=3;  H=4;  I=5;  J=6;  K=7;  L=8; M=9; N=10; O="Não" }
)

foreach ($rowData in $newRows) {
    Set-DataRow $ws $rowData["Row"] $rowData
}
